$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JAN-22")

# Fill in row 38 values (matches style/format of neighboring rows)
$ws.Range("A38").Value = 18
$ws.Range("B38").Value = 44592
$ws.Range("C38").Value = "RPA GSS"
$ws.Range("D38").Value = "1. Implementation of Public holidays at GSPN task is work in progress"
$ws.Range("E38").Value = 0.2
$ws.Range("F38").Value = "WIP"

# Copy number formats from similarly-formatted cells above so the same
# shared cellXf/style indices get reused instead of creating new ones.
$ws.Range("B36").Copy()
$ws.Range("B38").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E36").Copy()
$ws.Range("E38").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
